$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Addr, $Val)
    $c = $Sheet.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "62.260.93"
Set-TextCell $ws "E2" "  -0.03%  "
Set-TextCell $ws "D3" "3.027.15"
Set-TextCell $ws "E3" "  -0.22%  "
Set-TextCell $ws "E4" "  +0.03%  "
Set-TextCell $ws "D5" "541.39"
Set-TextCell $ws "E5" "  +0.94%  "
Set-TextCell $ws "D6" "133.41"
Set-TextCell $ws "E6" "  +0.41%  "
Set-TextCell $ws "D7" "1.00"
Set-TextCell $ws "E7" "  +0.09%  "
Set-TextCell $ws "D8" "3.019.31"
Set-TextCell $ws "E8" "  -0.32%  "
Set-TextCell $ws "D9" "0.493"
Set-TextCell $ws "E9" "  +0.79%  "
Set-TextCell $ws "D10" "6.13"
Set-TextCell $ws "E10" "  -0.73%  "
Set-TextCell $ws "E11" "  -3.51%  "
Set-TextCell $ws "E12" "  -0.21%  "
Set-TextCell $ws "D13" "0.0000222"
Set-TextCell $ws "E13" "  +1.37%  "
Set-TextCell $ws "D14" "34.24"
Set-TextCell $ws "E14" "  +1.15%  "
Set-TextCell $ws "D15" "3.513.05"
Set-TextCell $ws "E15" "  -0.26%  "
Set-TextCell $ws "D16" "62.238.52"
Set-TextCell $ws "E16" "  -0.09%  "
Set-TextCell $ws "D17" "3.023.82"
Set-TextCell $ws "E17" "  -0.25%  "
Set-TextCell $ws "E18" "  -3.37%  "
Set-TextCell $ws "D19" "6.62"
Set-TextCell $ws "E19" "  +1.12%  "
Set-TextCell $ws "D20" "479.63"
Set-TextCell $ws "E20" "  +3.71%  "
Set-TextCell $ws "D21" "13.27"
Set-TextCell $ws "E21" "  +0.18%  "
Set-TextCell $ws "E22" "  -1.73%  "
Set-TextCell $ws "E23" "  +1.78%  "
Set-TextCell $ws "D24" "80.80"
Set-TextCell $ws "E24" "  +4.04%  "
Set-TextCell $ws "D25" "12.10"
Set-TextCell $ws "E25" "  +1.29%  "
Set-TextCell $ws "E26" "  +0.10%  "
Set-TextCell $ws "E27" "  +1.19%  "
Set-TextCell $ws "D28" "7.74"
Set-TextCell $ws "E28" "  +0.18%  "
Set-TextCell $ws "D29" "1.00"
Set-TextCell $ws "E29" "  +0.06%  "
Set-TextCell $ws "E30" "  +4.22%  "
Set-TextCell $ws "D31" "25.66"
Set-TextCell $ws "E31" "  -0.16%  "
Set-TextCell $ws "E32" "  -1.29%  "
Set-TextCell $ws "D33" "5.65"
Set-TextCell $ws "E33" "  +4.91%  "
Set-TextCell $ws "D34" "2.36"
Set-TextCell $ws "E34" "  +4.17%  "
Set-TextCell $ws "D35" "55.00"
Set-TextCell $ws "E35" "  -5.43%  "
Set-TextCell $ws "D36" "5.87"
Set-TextCell $ws "E36" "  -0.21%  "
Set-TextCell $ws "D37" "459.55"
Set-TextCell $ws "E37" "  -0.35%  "
Set-TextCell $ws "D38" "3.166.09"
Set-TextCell $ws "E38" "  -0.55%  "
Set-TextCell $ws "D39" "0.0800"
Set-TextCell $ws "E39" "  +1.67%  "
Set-TextCell $ws "D40" "0.0387"
Set-TextCell $ws "E40" "  -0.27%  "
Set-TextCell $ws "E41" "  +1.58%  "
Set-TextCell $ws "D42" "8.08"
Set-TextCell $ws "E42" "  +0.68%  "
Set-TextCell $ws "D43" "2.47"
Set-TextCell $ws "E43" "  -0.54%  "
Set-TextCell $ws "D44" "26.40"
Set-TextCell $ws "E44" "  +5.93%  "
Set-TextCell $ws "D46" "0.244"
Set-TextCell $ws "E46" "  -0.67%  "
Set-TextCell $ws "E47" "  +0.94%  "
Set-TextCell $ws "D48" "1.97"
Set-TextCell $ws "E48" "  +0.45%  "
Set-TextCell $ws "B49" "PEPE"
Set-TextCell $ws "C49" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws "D49" "0.0₃0500"
Set-TextCell $ws "E49" "  -1.97%  "
Set-TextCell $ws "B50" "Monero"
Set-TextCell $ws "C50" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D50" "113.96"
Set-TextCell $ws "E50" "  -5.97%  "
Set-TextCell $ws "E51" "  +3.34%  "
